$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing the existing rows 14..125 down to 15..126.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new price-quote record (columns A,B,C,E,F,G,H,Q,R
# are identical to every other row in this block, and Excel's Insert already shifted
# the surrounding rows, so we only need to fill the record-specific columns here).
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 44490
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 100112028
$ws.Range("G14").Value = "Sandia"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 1000
$ws.Range("N14").Value = "$/kilo (volumen en unidades)"
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 1000
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
